$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(1,1).Select()
